# Update TPM-derived specificity/weight values on the active sheet
# (LR-pairs table) to reflect the newly recomputed TPM figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = 0.5597564161496535
$ws.Range("P2").Value = 0.5597564161496534
$ws.Range("S2").Value = 0.1074837636354556
$ws.Range("T2").Value = 0.1074837636354556

# Row 3
$ws.Range("M3").Value = 0.4111863333333334
$ws.Range("N3").Value = 1.233559
$ws.Range("O3").Value = 0.4402435838503465
$ws.Range("P3").Value = 0.4402435838503465
$ws.Range("Q3").Value = 0.09255681007433333
$ws.Range("R3").Value = 0.833011290669
$ws.Range("S3").Value = 0.08453505121760954
$ws.Range("T3").Value = 0.08453505121760951

# Row 4
$ws.Range("G4").Value = 0.9471683333333334
$ws.Range("I4").Value = 0.807981185146935
$ws.Range("O4").Value = 0.5597564161496535
$ws.Range("P4").Value = 0.5597564161496534
$ws.Range("R4").Value = 4.456710211665
$ws.Range("S4").Value = 0.452272652514198
$ws.Range("T4").Value = 0.4522726525141978

# Row 5
$ws.Range("G5").Value = 0.9471683333333334
$ws.Range("I5").Value = 0.807981185146935
$ws.Range("M5").Value = 0.4111863333333334
$ws.Range("N5").Value = 1.233559
$ws.Range("O5").Value = 0.4402435838503465
$ws.Range("P5").Value = 0.4402435838503465
$ws.Range("Q5").Value = 0.3894626740327778
$ws.Range("R5").Value = 3.505164066295
$ws.Range("S5").Value = 0.355708532632737
$ws.Range("T5").Value = 0.3557085326327369
